$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the A5 smoke-test label cell
$ws.Range("A5").Value = "103_TruckInsurance_001_SmokeTest_FillPage"

# Autofit column A to the new content
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update the active selection to A5
$ws.Range("A5").Select() | Out-Null
